# Add new "Carmel Cliff" community rows (lot 7362, four plans) to Sheet1's data table.
# Source data: community | lot | plan | work type | amount

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @"
1480|Carmel Cliff|7362|5528-2 Casena|FG|378
1481|Carmel Cliff|7362|5528-2 Casena|RG|953
1482|Carmel Cliff|7362|5528-2 Casena|LS|4255
1483|Carmel Cliff|7362|5528-2 Casena|PV|3151
1484|Carmel Cliff|7362|5528-2 Casena|PVO-70410|50
1485|Carmel Cliff|7362|5528-2 Casena|PVO-70468|0
1486|Carmel Cliff|7362|5528-2 Casena|PVO-70469|99
1487|Carmel Cliff|7362|5528-2 Casena|70706-LS|-83
1488|Carmel Cliff|7362|5528-2 Casena|70706-PV|353
1489|Carmel Cliff|7362|5528-2 Casena|PVO-70709|849
1490|Carmel Cliff|7362|5528-2 Casena|PVO-70714|1004
1491|Carmel Cliff|7362|5528-2 Casena|LSO-7440|0.75
1492|Carmel Cliff|7362|5528-2 Casena|LSO-74408|381
1493|Carmel Cliff|7362|5536-2-Matera|FG|378
1494|Carmel Cliff|7362|5536-2-Matera|RG|953
1495|Carmel Cliff|7362|5536-2-Matera|LS|4209
1496|Carmel Cliff|7362|5536-2-Matera|PV|3470
1497|Carmel Cliff|7362|5536-2-Matera|PVO-70410|50
1498|Carmel Cliff|7362|5536-2-Matera|PVO-70468|99
1499|Carmel Cliff|7362|5536-2-Matera|PVO-70469|0
1500|Carmel Cliff|7362|5536-2-Matera|70706-LS|0
1501|Carmel Cliff|7362|5536-2-Matera|70706-PV|0
1502|Carmel Cliff|7362|5536-2-Matera|PVO-70709|1655
1503|Carmel Cliff|7362|5536-2-Matera|PVO-70714|0
1504|Carmel Cliff|7362|5536-2-Matera|LSO-7440|0.75
1505|Carmel Cliff|7362|5536-2-Matera|LSO-74408|427
1506|Carmel Cliff|7362|5539-2-Pesaro|FG|378
1507|Carmel Cliff|7362|5539-2-Pesaro|RG|953
1508|Carmel Cliff|7362|5539-2-Pesaro|LS|4164
1509|Carmel Cliff|7362|5539-2-Pesaro|PV|3562
1510|Carmel Cliff|7362|5539-2-Pesaro|PVO-70410|50
1511|Carmel Cliff|7362|5539-2-Pesaro|PVO-70468|0
1512|Carmel Cliff|7362|5539-2-Pesaro|PVO-70469|99
1513|Carmel Cliff|7362|5539-2-Pesaro|70706-LS|0
1514|Carmel Cliff|7362|5539-2-Pesaro|70706-PV|0
1515|Carmel Cliff|7362|5539-2-Pesaro|PVO-70709|1390
1516|Carmel Cliff|7362|5539-2-Pesaro|PVO-70714|0
1517|Carmel Cliff|7362|5539-2-Pesaro|LSO-7440|0.75
1518|Carmel Cliff|7362|5539-2-Pesaro|LSO-74408|472
1519|Carmel Cliff|7362|5545-1-Victoria|FG|378
1520|Carmel Cliff|7362|5545-1-Victoria|RG|953
1521|Carmel Cliff|7362|5545-1-Victoria|LS|4164
1522|Carmel Cliff|7362|5545-1-Victoria|PV|3541
1523|Carmel Cliff|7362|5545-1-Victoria|PVO-70410|50
1524|Carmel Cliff|7362|5545-1-Victoria|PVO-70468|0
1525|Carmel Cliff|7362|5545-1-Victoria|PVO-70469|99
1526|Carmel Cliff|7362|5545-1-Victoria|70706-LS|0
1527|Carmel Cliff|7362|5545-1-Victoria|70706-PV|0
1528|Carmel Cliff|7362|5545-1-Victoria|PVO-70709|1379
1529|Carmel Cliff|7362|5545-1-Victoria|PVO-70714|1655
1530|Carmel Cliff|7362|5545-1-Victoria|LSO-7440|0.75
1531|Carmel Cliff|7362|5545-1-Victoria|LSO-74408|472
"@

$lines = $newRows -split "`r?`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $lines) {
    $f = $line -split '\|'
    $r = [int]$f[0]
    $community = $f[1]
    $lot = [double]$f[2]
    $plan = $f[3]
    $workType = $f[4]
    $amount = [double]$f[5]

    $ws.Cells.Item($r, 1).Value = $community
    $ws.Cells.Item($r, 2).Value = $lot
    $ws.Cells.Item($r, 3).Value = $plan
    $ws.Cells.Item($r, 4).Value = $workType
    $ws.Cells.Item($r, 5).Value = $amount
}

$lastRow = 1531

# Resize the data table (ListObject) to include the newly appended rows, which also
# extends the AutoFilter range to match.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:E" + $lastRow))

# Restore the view/selection state recorded in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J57").Select()
